$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Headers: BTec_Logo-Orange (jpg) — swap name2 -> name1
$h1 = $sec.Headers.Item(1)
$h1.Range.InlineShapes.Item(1).Name = "image1.jpg"

$h2 = $sec.Headers.Item(2)
$h2.Range.InlineShapes.Item(1).Name = "image1.jpg"

# Footers: Pearson logo (png) — swap name1 -> name2
$f1 = $sec.Footers.Item(1)
$f1.Range.InlineShapes.Item(1).Name = "image2.png"

$f2 = $sec.Footers.Item(2)
$f2.Range.InlineShapes.Item(1).Name = "image2.png"
